$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A24").Value = 23
$ws.Range("B24").Value = "Estatística"
$ws.Range("C24").Value = "Medidas de Variabilidade"
$ws.Range("D24").Value = "Fórmula do <b>desvio padrão</b>`n<i>amostral e populacional</i>"
$ws.Range("E24").Value = "populacional: sig<sup>2</sup> = (sum<sub>i=1</sub><sup>n</sup>(x<sub>i</sub> - mu)<sup>2</sup>)/n`namostral: S<sup>2</sup> = (sum<sub>i=1</sub><sup>n</sup>(x<sub>i</sub> - x<sup>-</sup>)<sup>2</sup>)/(n-1)"
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 0

$ws.Range("A25").Value = 24
$ws.Range("B25").Value = "Conhecimentos Específicos"
$ws.Range("C25").Value = "Gestão da Manutenção e Confiabilidade"
$ws.Range("D25").Value = "Tipos de manutenção"
$ws.Range("E25").Value = "<b>Manutenção não-planejada</b>`n<ul>`n	<li>Corretiva</li>`n</ul>`n<b>Manutenção planejada</b>`n<ul>`n	<li>Corretiva</li>`n	<li>Preventiva</li>`n	<li>Preditiva</li>`n	<li>Prescritiva</li>`n	<li>Detectiva</li>`n</ul>"
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0

$ws.Range("A26").Value = 25
$ws.Range("B26").Value = "Conhecimentos Específicos"
$ws.Range("C26").Value = "Gestão da Manutenção e Confiabilidade"
$ws.Range("D26").Value = "Princípios Fundamentais da TPM"
$ws.Range("E26").Value = "<ul>`n	<li>Maximizar a eficiência global</li>`n	<li>Envolvimento de todos os níveis</li>`n	<li>Prevenção de falhas</li>`n	<li>Abordagem sistemática</li>`n</ul>"
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = 0

$ws.Range("A27").Value = 26
$ws.Range("B27").Value = "Conhecimentos Específicos"
$ws.Range("C27").Value = "Gestão da Manutenção e Confiabilidade"
$ws.Range("D27").Value = "Objetivos da TPM"
$ws.Range("E27").Value = "<ul>`n	<li>Aumentar a disponibilidade</li>`n	<li>Reduzir custos operacionais</li>`n	<li>Melhorar a segurança</li>`n	<li>Melhorar a segurança</li>`n</ul>"
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 0

$ws.Range("A28").Value = 27
$ws.Range("B28").Value = "Conhecimentos Específicos"
$ws.Range("C28").Value = "Gestão da Manutenção e Confiabilidade"
$ws.Range("D28").Value = "Os 8 Pilares da TPM"
$ws.Range("E28").Value = "<ul>`n	<li>Manutenção autônoma</li>`n	<li>Manutenção planejada</li>`n	<li>Educação e treinamento</li>`n	<li>Melhoria focada</li>`n	<li>Gestão da segurança e do meio ambiente</li>`n	<li>Manutenção de qualidade</li>`n	<li>Controle inicial</li>`n	<li>Gestão administrativa</li>`n</ul>"
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 0

$ws.Range("A29").Value = 28
$ws.Range("B29").Value = "Conhecimentos Específicos"
$ws.Range("C29").Value = "Gestão da Manutenção e Confiabilidade"
$ws.Range("D29").Value = "seis grandes perdas da TPM"
$ws.Range("E29").Value = "<ul>`n	<li>por quebra de equipamento</li>`n	<li>decorrentes de ajustes nas preparações</li>`n	<li>nas paradas curtas e frequentes</li>`n	<li>por uma operação abaixo da normal</li>`n	<li>decorrentes de peças defeituosas e retrabalhos</li>`n	<li>provenientes do início da produção</li>`n</ul>"
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 0

$ws.Range("A30").Value = 29
$ws.Range("B30").Value = "Conhecimentos Específicos"
$ws.Range("C30").Value = "Gestão da Manutenção e Confiabilidade"
$ws.Range("D30").Value = "três fatores importantes da TPM"
$ws.Range("E30").Value = "<ul>`n	<li>busca pelo lucro através da economicidade</li>`n	<li>ser um sistema integrado</li>`n	<li>próprio operador executa a manutenção</li>`n</ul>"
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0

$ws.Range("A31").Value = 30
$ws.Range("B31").Value = "Conhecimentos Específicos"
$ws.Range("C31").Value = "Gestão da Manutenção e Confiabilidade"
$ws.Range("D31").Value = "Pontos-chave para implementação da TPM:"
$ws.Range("E31").Value = "<ul>`n	<li>Capacitação;</li>`n	<li>Aplicar o programa 5s/8s;</li>`n	<li>Eliminar as 6 grandes perdas;</li>`n	<li>Aplicar as 5 ações para alcançar a `"quebra zero`":</li>`n</ul>"
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 0

$ws.Range("A32").Value = "'31"
$ws.Range("B32").Value = "Conhecimentos Específicos"
$ws.Range("C32").Value = "Gestão da Manutenção e Confiabilidade"
$ws.Range("D32").Value = "Fórmula do MTBF"
$ws.Range("E32").Value = "<b>Tempo Médio Entre Falhas</b>`nMTBF = (Tempo Total Disponível - Tempo das Paradas ou Tempo Perdido)/Quantidade de Paradas;"
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 0

$ws.Range("A33").Value = 32
$ws.Range("B33").Value = "Conhecimentos Específicos"
$ws.Range("C33").Value = "Gestão da Manutenção e Confiabilidade"
$ws.Range("D33").Value = "Fórmula do MTTR"
$ws.Range("E33").Value = "<b>Tempo Médio de Reparo</b>`nMTTR = (Tempo de Parada)/Quantidade de Paradas;"
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 0

$ws.Range("A34").Value = 33
$ws.Range("B34").Value = "Conhecimentos Específicos"
$ws.Range("C34").Value = "Gestão da Manutenção e Confiabilidade"
$ws.Range("D34").Value = "<b>Disponibilidade</b>`n<i>Definição e Fórmula</i>"
$ws.Range("E34").Value = "<b>Definição</b>`nexpressa como a porcentagem de tempo em que um ativo está operando, em comparação com o tempo total de operação programado`nDisponibilidade = MTBF/(MTBF + MTTR)`n<b>ou</b>`nDisponibilidade = Tempo disponível/(Tempo disponível + Tempo em manutenção)"
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 0

$ws.Range("A35").Value = 34
$ws.Range("B35").Value = "Conhecimentos Específicos"
$ws.Range("C35").Value = "Gestão da Manutenção e Confiabilidade"
$ws.Range("D35").Value = "<b>Taxa de Falha (λ)</b>"
$ws.Range("E35").Value = "Taxa de Falha (λ) = 1/MTBF"
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 0

$ws.Range("A36").Value = 35
$ws.Range("B36").Value = "Conhecimentos Específicos"
$ws.Range("C36").Value = "Gestão da Manutenção e Confiabilidade"
$ws.Range("D36").Value = "Manutenabilidade"
$ws.Range("E36").Value = "facilidade com que as atividades de manutenção podem ser realizadas em um ativo ou equipamento"
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 0

$ws.Range("A37").Value = 36
$ws.Range("B37").Value = "Conhecimentos Específicos"
$ws.Range("C37").Value = "Gestão da Manutenção e Confiabilidade"
$ws.Range("D37").Value = "Confiabilidade"
$ws.Range("E37").Value = "probabilidade de que um produto, sistema ou serviço desempenhe sua função pretendida adequadamente por um período de tempo especificado"
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 0

$ws.Range("A38").Value = 37
$ws.Range("B38").Value = "Conhecimentos Específicos"
$ws.Range("C38").Value = "Gestão da Manutenção e Confiabilidade"
$ws.Range("D38").Value = "Curva P-F"
$ws.Range("E38").Value = "Curva que mostra a performance do equipamento em função do tempo, trazendo dois pontos.`n<ul>`n	<li>P - Falha Pontencial: Momento em que é identificado que o equipamento está prestes a falhar</li>`n	<li>F - Falha Funcional: Momento em que o equipamento falha</li>`n</ul>"
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 0

$ws.Range("A39").Value = 38
$ws.Range("B39").Value = "Conhecimentos Específicos"
$ws.Range("C39").Value = "Gestão da Manutenção e Confiabilidade"
$ws.Range("D39").Value = "Fórmula da Confiabilidade"
$ws.Range("E39").Value = "R(t) = e<sup>-λt</sup>"
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 0

$ws.Range("A40").Value = 39
$ws.Range("B40").Value = "Conhecimentos Específicos"
$ws.Range("C40").Value = "Gestão da Manutenção e Confiabilidade"
$ws.Range("D40").Value = "Formulas para a confiabilidade de sistemas <b>em série</b> e <b>em paralelo</b>"
$ws.Range("E40").Value = "<b>Em série</b>`nR<sub>s</sub> = P<sub>1</sub> × P<sub>2</sub> × ... × P<sub>n</sub>`n<b>Em paralelo</b>`nRs = 1 - [(1 - P<sub>1</sub>) x (1 - P<sub>2</sub>) x (1 - P<sub>3</sub>) x ... x (1 - P<sub>n</sub>)]"
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 0

$ws.Range("A41").Value = 40
$ws.Range("B41").Value = "Conhecimentos Específicos"
$ws.Range("C41").Value = "Gestão da Manutenção e Confiabilidade"
$ws.Range("D41").Value = "etapas do FMEA"
$ws.Range("E41").Value = "<ol>`n	<li>Função do equipamento</li>`n	<li>Falha funcional</li>`n	<li>Componente</li>`n	<li>Modo (s) de falha potencial</li>`n	<li>Efeito (s) potencial (is) da falha</li>`n	<li>Causa (s) potencial (is) de falha</li>`n	<li>Controles atuais</li>`n	<li>Frequência do processo atual</li>`n</ol>"
$ws.Range("F41").Value = 0
$ws.Range("G41").Value = 0

$ws.Range("A42").Value = 41
$ws.Range("B42").Value = "Conhecimentos Específicos"
$ws.Range("C42").Value = "Gestão da Manutenção e Confiabilidade"
$ws.Range("D42").Value = "Etapas do processo do RCFA"
$ws.Range("E42").Value = "<ol>`n	<li>Identificação dos sintomas</li>`n	<li>Avaliação das causas endereçáveis</li>`n	<li>Coletando e analisando dados</li>`n	<li>Isolando e testando variáveis</li>`n	<li>Identificando a (s) causa (s)</li>`n	<li>Criação e implementação de um plano de ação</li>`n</ol>"
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 0

$ws.Range("A43").Value = 42
$ws.Range("B43").Value = "Conhecimentos Específicos"
$ws.Range("C43").Value = "Gestão da Manutenção e Confiabilidade"
$ws.Range("D43").Value = "3 fases da FTA"
$ws.Range("E43").Value = "<ol>`n	<li>Identificar o perigo;</li>`n	<li>Obter entendimento do sistema que está sendo analisado;</li>`n	<li>Criar a árvore de falhas;</li>`n</ol>"
$ws.Range("F43").Value = 0
$ws.Range("G43").Value = 0

$ws.Range("A44").Value = 43
$ws.Range("B44").Value = "Conhecimentos Específicos"
$ws.Range("C44").Value = "Gestão da Manutenção e Confiabilidade"
$ws.Range("D44").Value = "Fórmula da função densidade de probabilidade exponencial f(t)"
$ws.Range("E44").Value = "f(t) = λe<sup>−λt</sup>, t ≥ 0"
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 0

$ws.Range("A45").Value = 44
$ws.Range("B45").Value = "Conhecimentos Específicos"
$ws.Range("C45").Value = "Gestão da Manutenção e Confiabilidade"
$ws.Range("D45").Value = "Fórmula da função de distribuição exponencial acumulada F(t)"
$ws.Range("E45").Value = "F(t) = 1−e<sup>−λt</sup>"
$ws.Range("F45").Value = 0
$ws.Range("G45").Value = 0

$ws.Range("A46").Value = 45
$ws.Range("B46").Value = "Conhecimentos Específicos"
$ws.Range("C46").Value = "Gestão da Manutenção e Confiabilidade"
$ws.Range("D46").Value = "<b>Função confiabilidade</b>`n<i>Definição e fórmula</i>"
$ws.Range("E46").Value = "<b>Definição</b>`nrepresenta a probabilidade de que um equipamento permaneça em funcionamento até o tempo t, <u>sem apresentar falhas</u>:`n<b>Fórmula</b>`nR(t) = e<sup>−λt</sup>"
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 0

$ws.Range("A47").Value = 46
$ws.Range("B47").Value = "Conhecimentos Específicos"
$ws.Range("C47").Value = "Gestão da Manutenção e Confiabilidade"
$ws.Range("D47").Value = "Fórmula para estimar a taxa de falhas"
$ws.Range("E47").Value = "λ=k/n⋅Δt `n<ul>`n	<li><i>onde:</i></li>`n	<li>λ: taxa de falha estimada;</li>`n	<li>k: número total de falhas registradas;</li>`n	<li>n: número de equipamentos monitorados (ou unidades observadas);</li>`n	<li>Δt: intervalo de tempo de observação (em horas, dias, semanas etc.).</li>`n</ul>"
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 0
